# Commit: "Added one more word in matrix generator comment"
# The comment cell (H21) that explains the random-matrix generator gets the
# word "файла" appended, so the sentence reads "...при любом изменении файла)"
# instead of "...при любом изменении)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H21").Value = "СЛУЧАЙНАЯ МАТРИЦА (генерируется заново при любом изменении файла)"

# The author's last selection in the saved file was this same cell.
$ws.Range("H21").Select() | Out-Null
